$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text first so numeric-looking strings
# (e.g. "568.90", "1.00") are stored as text, not auto-converted to numbers.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '64.366.01'
$ws.Range('E2').Value = '  +0.10%  '

$ws.Range('D3').Value = '3.111.58'
$ws.Range('E3').Value = '  -1.58%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '568.90'
$ws.Range('E5').Value = '  -0.31%  '

$ws.Range('D6').Value = '161.47'
$ws.Range('E6').Value = '  -3.99%  '

$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  -5.29%  '

$ws.Range('D9').Value = '3.123.54'
$ws.Range('E9').Value = '  -1.93%  '

$ws.Range('E10').Value = '  -1.68%  '

$ws.Range('D11').Value = '6.61'
$ws.Range('E11').Value = '  -3.08%  '

$ws.Range('D12').Value = '0.379'
$ws.Range('E12').Value = '  -2.47%  '

$ws.Range('D13').Value = '3.666.54'
$ws.Range('E13').Value = '  -1.18%  '

$ws.Range('E14').Value = '  -2.30%  '

$ws.Range('D15').Value = '64.465.47'
$ws.Range('E15').Value = '  +0.11%  '

$ws.Range('D16').Value = '24.65'
$ws.Range('E16').Value = '  -2.63%  '

$ws.Range('D17').Value = '3.122.39'
$ws.Range('E17').Value = '  -1.15%  '

$ws.Range('E18').Value = '  -1.40%  '

$ws.Range('D19').Value = '405.86'
$ws.Range('E19').Value = '  -2.80%  '

$ws.Range('D20').Value = '5.19'
$ws.Range('E20').Value = '  -1.98%  '

$ws.Range('D21').Value = '12.33'
$ws.Range('E21').Value = '  -3.96%  '

$ws.Range('D22').Value = '6.97'
$ws.Range('E22').Value = '  -2.39%  '

$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.19%  '

$ws.Range('D24').Value = '67.74'
$ws.Range('E24').Value = '  -2.61%  '

$ws.Range('D25').Value = '0.478'
$ws.Range('E25').Value = '  -3.84%  '

$ws.Range('E26').Value = '  -4.88%  '

$ws.Range('D27').Value = '0.0000102'
$ws.Range('E27').Value = '  -1.28%  '

$ws.Range('D28').Value = '9.08'
$ws.Range('E28').Value = '  +3.50%  '

$ws.Range('D29').Value = '0.998'

$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.12%  '

$ws.Range('D31').Value = '1.79'
$ws.Range('E31').Value = '  -1.62%  '

$ws.Range('D32').Value = '21.18'
$ws.Range('E32').Value = '  -2.48%  '

$ws.Range('D33').Value = '164.62'
$ws.Range('E33').Value = '  +6.10%  '

$ws.Range('D34').Value = '4.89'
$ws.Range('E34').Value = '  -3.36%  '

$ws.Range('D35').Value = '6.21'
$ws.Range('E35').Value = '  -2.26%  '

$ws.Range('D36').Value = '1.12'
$ws.Range('E36').Value = '  +0.41%  '

$ws.Range('E37').Value = '  -0.96%  '

$ws.Range('E38').Value = '  -2.34%  '

$ws.Range('D39').Value = '2.594.24'
$ws.Range('E39').Value = '  -4.02%  '

$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '23.67'
$ws.Range('E40').Value = '  -0.89%  '

$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.10'
$ws.Range('E41').Value = '  -2.65%  '

$ws.Range('E42').Value = '  -1.83%  '

$ws.Range('D43').Value = '0.687'
$ws.Range('E43').Value = '  -4.21%  '

$ws.Range('D44').Value = '0.0615'
$ws.Range('E44').Value = '  -0.39%  '

$ws.Range('D45').Value = '5.20'
$ws.Range('E45').Value = '  -4.94%  '

$ws.Range('D46').Value = '0.0253'
$ws.Range('E46').Value = '  -3.59%  '

$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '20.92'
$ws.Range('E47').Value = '  -1.94%  '

$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '283.74'
$ws.Range('E48').Value = '  -1.71%  '

$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.05%  '

$ws.Range('E50').Value = '  -1.77%  '

$ws.Range('E51').Value = '  +0.29%  '

# Restore the default (Normal) style so no stray number-format/quote-prefix
# style survives on cells that must look identical to the original sheet.
$ws.Range('D2:E51').Style = 'Normal'
